$d = $word.ActiveDocument

$replacements = @(
    @("2025-11-27 Thursday", "2025-11-28 Friday"),
    @("85×80=", "65×12="),
    @("12×98=", "86×48="),
    @("33×75=", "64×60="),
    @("87×47=", "37×65="),
    @("88×14=", "57×40="),
    @("51×93=", "83×85="),
    @("49×70=", "33×73="),
    @("18×16=", "14×80="),
    @("69×78=", "91×97="),
    @("43×70=", "36×62="),
    @("81×60=", "47×21="),
    @("73×13=", "27×92="),
    @("40×72=", "60×96="),
    @("68×11=", "32×84="),
    @("88×20=", "26×63="),
    @("72×44=", "18×58="),
    @("45×28=", "63×87="),
    @("47×49=", "26×94="),
    @("44×64=", "64×97="),
    @("12×92=", "67×26="),
    @("32×98=", "79×39="),
    @("87×25=", "84×74="),
    @("72×22=", "38×65="),
    @("63×57=", "27×67="),
    @("17×56=", "25×29=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
